$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 12 new rows (148-159) to the "C00 / CCP 장비코드" common-code category,
# cloning the format (font/style + row height) of the last existing data row (147)
# and then filling in the real values.
$srcRow = $ws.Range("A147:E147")
$srcRow.Copy($ws.Range("A148:E148"))
$ws.Rows.Item(148).RowHeight = 15.75
$srcRow.Copy($ws.Range("A149:E149"))
$ws.Rows.Item(149).RowHeight = 15.75
$srcRow.Copy($ws.Range("A150:E150"))
$ws.Rows.Item(150).RowHeight = 15.75
$srcRow.Copy($ws.Range("A151:E151"))
$ws.Rows.Item(151).RowHeight = 15.75
$srcRow.Copy($ws.Range("A152:E152"))
$ws.Rows.Item(152).RowHeight = 15.75
$srcRow.Copy($ws.Range("A153:E153"))
$ws.Rows.Item(153).RowHeight = 15.75
$srcRow.Copy($ws.Range("A154:E154"))
$ws.Rows.Item(154).RowHeight = 15.75
$srcRow.Copy($ws.Range("A155:E155"))
$ws.Rows.Item(155).RowHeight = 15.75
$srcRow.Copy($ws.Range("A156:E156"))
$ws.Rows.Item(156).RowHeight = 15.75
$srcRow.Copy($ws.Range("A157:E157"))
$ws.Rows.Item(157).RowHeight = 15.75
$srcRow.Copy($ws.Range("A158:E158"))
$ws.Rows.Item(158).RowHeight = 15.75
$srcRow.Copy($ws.Range("A159:E159"))
$ws.Rows.Item(159).RowHeight = 15.75

$ws.Cells.Item(148,1).Value = "C00"
$ws.Cells.Item(148,2).Value = "`$`$"
$ws.Cells.Item(148,3).Value = "CCP 장비코드"
$ws.Cells.Item(148,5).Value = 20201201153351

$ws.Cells.Item(149,1).Value = "C00"
$ws.Cells.Item(149,2).Value = "CONOVN"
$ws.Cells.Item(149,3).Value = "컨벤션오븐기"
$ws.Cells.Item(149,5).Value = 20201201153352

$ws.Cells.Item(150,1).Value = "C00"
$ws.Cells.Item(150,2).Value = "FREZFS"
$ws.Cells.Item(150,3).Value = "급냉고(반제품)-2"
$ws.Cells.Item(150,5).Value = 20201201153353

$ws.Cells.Item(151,1).Value = "C00"
$ws.Cells.Item(151,2).Value = "FREZMA"
$ws.Cells.Item(151,3).Value = "냉동고(원재료)"
$ws.Cells.Item(151,5).Value = 20201201153354

$ws.Cells.Item(152,1).Value = "C00"
$ws.Cells.Item(152,2).Value = "QUFRE1"
$ws.Cells.Item(152,3).Value = "급냉고(반제품)-1"
$ws.Cells.Item(152,5).Value = 20201201153355

$ws.Cells.Item(153,1).Value = "C00"
$ws.Cells.Item(153,2).Value = "QUFRE2"
$ws.Cells.Item(153,3).Value = "냉동고(완제품)"
$ws.Cells.Item(153,5).Value = 20201201153356

$ws.Cells.Item(154,1).Value = "C00"
$ws.Cells.Item(154,2).Value = "BBMIX1"
$ws.Cells.Item(154,3).Value = "볶음솥+밥혼합기"
$ws.Cells.Item(154,5).Value = 20201201153357

$ws.Cells.Item(155,1).Value = "C00"
$ws.Cells.Item(155,2).Value = "TART1"
$ws.Cells.Item(155,3).Value = "타르트제조용"
$ws.Cells.Item(155,5).Value = 20201201153358

$ws.Cells.Item(156,1).Value = "C00"
$ws.Cells.Item(156,2).Value = "RTEMPS"
$ws.Cells.Item(156,3).Value = "상온저장고(원재료)"
$ws.Cells.Item(156,5).Value = 20201201153359

$ws.Cells.Item(157,1).Value = "C00"
$ws.Cells.Item(157,2).Value = "COOKRM"
$ws.Cells.Item(157,3).Value = "조리실 룸"
$ws.Cells.Item(157,5).Value = 20201201153360

$ws.Cells.Item(158,1).Value = "C00"
$ws.Cells.Item(158,2).Value = "COLDSM"
$ws.Cells.Item(158,3).Value = "저온저장고(원재료)"
$ws.Cells.Item(158,5).Value = 20201201153361

$ws.Cells.Item(159,1).Value = "C00"
$ws.Cells.Item(159,2).Value = "REFRCD"
$ws.Cells.Item(159,3).Value = "냉장고(식힘)"
$ws.Cells.Item(159,5).Value = 20201201153362

# Match the author's final scroll/selection state: cursor on the first new row,
# sheet scrolled down so that row 133 is at the top.
$ws.Range("B148").Select()
$excel.ActiveWindow.ScrollRow = 133